$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize / reposition the workbook window to match the saved view state.
$excel.ActiveWindow.Width = 16800
$excel.ActiveWindow.Height = 9680
$excel.ActiveWindow.Top = 1920
$excel.ActiveWindow.Left = 680

# Update the route-point coordinate columns (C and D) for rows 2-5.
$ws.Range("C2:C5").Value = 97.8473625
$ws.Range("D2:D5").Value = 2.8316146

# Move the active selection in the frozen (bottom-right) pane to D12.
$ws.Range("D12").Select()
